# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# on the cryptos worksheet, per commit "Updated cryptos list on Thu May 18
# 04:57:48 UTC 2023 with GitHub Actions".
#
# Column D holds numeric-looking price strings ("314.60", "1.200",
# "0.000008797", ...) that must remain plain TEXT, matching the original
# inlineStr cells, rather than being auto-converted to numbers (which would
# silently drop meaningful trailing/leading zeros). A handful of D values
# (multi-dot "thousands" prices like "27.342.44") are already unambiguous
# text to Excel and need no special handling; the rest are prefixed with a
# leading apostrophe - the standard Excel way to force text entry - so they
# are stored as text without touching the cell's number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.342.44"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").Value = "1.825.39"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("D5").Value = "'314.60"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'0.4484"
$ws.Range("E7").Value = "  -2.04%  "

$ws.Range("D8").Value = "'0.3779"
$ws.Range("E8").Value = "  +1.13%  "

$ws.Range("D9").Value = "'0.07440"
$ws.Range("E9").Value = "  +1.69%  "

$ws.Range("D10").Value = "'0.8881"
$ws.Range("E10").Value = "  +2.89%  "

$ws.Range("D11").Value = "'20.99"
$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("D12").Value = "1.827.39"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").Value = "'6.741"
$ws.Range("E13").Value = "  +0.35%  "

$ws.Range("D14").Value = "'5.463"
$ws.Range("E14").Value = "  +1.72%  "

$ws.Range("D15").Value = "'93.55"
$ws.Range("E15").Value = "  +0.63%  "

$ws.Range("D16").Value = "'0.07133"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "'0.000008797"
$ws.Range("E18").Value = "  -0.70%  "

$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").Value = "'15.14"
$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("D21").Value = "27.344.12"
$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("D22").Value = "'5.399"
$ws.Range("E22").Value = "  +3.75%  "

$ws.Range("D23").Value = "'10.98"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").Value = "2.050.04"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "'1.970"

$ws.Range("D26").Value = "'151.53"
$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("D27").Value = "'2.329"
$ws.Range("E27").Value = "  +4.58%  "

$ws.Range("E28").Value = "  +0.79%  "

$ws.Range("D29").Value = "'5.393"
$ws.Range("E29").Value = "  +1.98%  "

$ws.Range("D30").Value = "'117.81"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").Value = "'0.08892"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").Value = "'0.7939"
$ws.Range("E32").Value = "  +4.16%  "

$ws.Range("D33").Value = "'1.200"
$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "'4.609"
$ws.Range("E34").Value = "  +2.86%  "

$ws.Range("D35").Value = "'2.920"
$ws.Range("E35").Value = "  -1.84%  "

$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "'1.111"
$ws.Range("E37").Value = "  +0.50%  "

$ws.Range("D38").Value = "'0.01985"
$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("D39").Value = "'0.05314"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "'7.304"
$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("D41").Value = "'0.5345"
$ws.Range("E41").Value = "  -0.64%  "

$ws.Range("D42").Value = "'2.873"

$ws.Range("D43").Value = "'0.1717"
$ws.Range("E43").Value = "  -0.17%  "

$ws.Range("D44").Value = "'2.319"
$ws.Range("E44").Value = "  +16.73%  "

$ws.Range("D45").Value = "'8.672"
$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("D46").Value = "'0.5080"
$ws.Range("E46").Value = "  -2.77%  "

$ws.Range("D47").Value = "'10.66"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("D48").Value = "'1.697"
$ws.Range("E48").Value = "  +0.75%  "

$ws.Range("D49").Value = "'105.23"
$ws.Range("E49").Value = "  -0.80%  "

$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("D51").Value = "'0.06411"
$ws.Range("E51").Value = "  -0.18%  "
